$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Sph"
$ws.Range("C6").Value = 0.6035
$ws.Range("D6").Value = 1.3499
$ws.Range("E6").Value = 0.21
$ws.Range("F6").Value = 0.4470701533446922
$ws.Range("G6").Value = 840.9776000000001
$ws.Range("H6").Value = 0.8457204115375917
